$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing rows down
$ws.Rows.Item(1).Insert()

# Populate the new first row with the new words
$ws.Range("A1").Value = "exceção"
$ws.Range("B1").Value = "concessão"
$ws.Range("C1").Value = "impressão"
$ws.Range("D1").Value = "presunção"
$ws.Range("E1").Value = "concepção"
$ws.Range("F1").Value = "inspiração"
